$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    8 = @{ 4=22610000; 5=16299000; 6=16768900; 7=17887200; 8=17661100; 9=17040300; 10=17915100 }
    9 = @{ 4=15963100; 5=11629300; 6=11894600; 7=12666800; 8=12593200; 9=12452200; 10=13024500 }
    10 = @{ 4=6646900; 5=4669700; 6=4874300; 7=5220400; 8=5067900; 9=4588100; 10=4890600 }
    12 = @{ 4=665600; 5=637400; 6=639500; 7=639300; 8=582900; 9=549500 }
    14 = @{ 4=181200; 5=15800; 6=27400; 7=10200; 8=20800; 9=17200; 10=28100 }
    17 = @{ 4=20154900; 5=14725200; 6=14883300; 7=15699000; 8=15487000; 9=15127400; 10=15597800 }
    18 = @{ 4=2455100; 5=1573800; 6=1885500; 7=2188200; 8=2174100; 9=1912900; 10=2317300 }
    20 = @{ 4=348900; 5=5300; 6=45900; 7=30200; 8=93900; 9=11200; 10=9500 }
    21 = @{ 4=4015900; 5=2523000; 6=2949200; 7=3143700; 8=3057800; 9=2732400; 10=3142300 }
    22 = @{ 4=166100; 5=74200; 6=79300; 7=84300; 8=79800; 9=74500; 10=70400 }
    23 = @{ 4=2637900; 5=1504900; 6=1852100; 7=2134100; 8=2188200; 9=1849600; 10=2256500 }
    24 = @{ 4=780900; 5=455700; 6=576000; 7=709600; 8=686500; 9=624600; 10=673200 }
    26 = @{ 4=1857000; 5=1049200; 6=1276100; 7=1424500; 8=1501700; 9=1225000; 10=1583300 }
    27 = @{ 4=1775500; 5=1025000; 6=1242300; 7=1392200; 8=1442000; 9=1141900; 10=1510100 }
    32 = @{ 4=-348900; 5=-5300; 6=-45900; 7=-30200; 8=-93900; 9=-11200; 10=-9500 }
    33 = @{ 4=1775500; 5=1025000; 6=1242300; 7=1392200; 8=1442000; 9=1141900; 10=1510100 }
    35 = @{ 4=1775500; 5=1025000; 6=1242300; 7=1392200; 8=1442000; 9=1141900; 10=1510100 }
    41 = @{ 4=1327600; 5=1104600; 6=980600; 7=970100; 8=824000; 9=848300; 10=751000 }
    42 = @{ 4=9000; 6=2300; 7=3400; 8=4400; 9=5100; 10=8200 }
    43 = @{ 4=7166200; 5=5598300; 6=5274000; 7=5605600; 8=5580900; 9=5486700; 10=5060100 }
    44 = @{ 4=6601800; 5=4826400; 6=4878100; 7=5630800; 8=5650700; 9=5728200; 10=5535700 }
    45 = @{ 4=1145600; 5=1301500; 6=1398600; 7=1543800; 8=1440700; 9=1419900; 10=1304300 }
    46 = @{ 4=16250200; 5=12832500; 6=12533500; 7=13753800; 8=13500600; 9=13488200; 10=12659400 }
    47 = @{ 4=3656300; 5=3724400; 6=3359600; 7=3452500; 8=3175500; 9=2843200; 10=2365200 }
    48 = @{ 4=6694400; 5=6138400; 6=6307600; 7=6725000; 8=6032800; 9=5290400; 10=4788100 }
    49 = @{ 4=2975000; 5=914400; 6=931700; 7=858600; 8=854800; 9=842800; 10=806200 }
    52 = @{ 4=911800; 5=404900; 6=504100; 7=507700; 8=406400; 9=296900; 10=358800 }
    54 = @{ 4=30487700; 5=24014600; 6=23636500; 7=25297600; 8=23970100; 9=22761400; 10=20977600 }
    57 = @{ 4=2744100; 5=2170600; 6=1856900; 7=2034800; 8=2117400; 9=2045500; 10=2472100 }
    58 = @{ 4=2982300; 5=1969300; 6=2214000; 7=2801100; 8=2657800; 9=3037000; 10=3030900 }
    59 = @{ 4=3220100; 5=2189700; 6=2265100; 7=2440400; 8=2386600; 9=2398800; 10=2304900 }
    60 = @{ 4=8946500; 5=6329600; 6=6336100; 7=7276400; 8=7161800; 9=7481300; 10=7807900 }
    61 = @{ 4=4345500; 5=1725400; 6=1922200; 7=2524600; 8=2812000; 9=3108100; 10=2825200 }
    62 = @{ 4=1433600; 5=1057000; 6=1024800; 7=1046200; 8=968600; 9=847700; 10=785100 }
    66 = @{ 4=15440300; 5=9761500; 6=9919000; 7=11475700; 8=11527500; 9=11975000; 10=11849900 }
    72 = @{ 4=13901600; 5=12680600; 6=12150200; 7=11402300; 8=10682700; 9=9697500; 10=8943700 }
    76 = @{ 4=15047400; 5=14253100; 6=13717400; 7=13821900; 8=12442600; 9=10786500; 10=9127700 }
    81 = @{ 4=1775500; 5=1025000; 6=1242300; 7=1392200; 8=1442000; 9=1141900; 10=1510100 }
    83 = @{ 4=1210500; 5=942800; 6=1016700; 7=924200; 8=788900; 9=807500; 10=814600 }
    89 = @{ 4=1341500; 5=2315400; 6=2889500; 7=3106600; 8=2887600; 9=1935000; 10=954700 }
    91 = @{ 4=-1400500; 5=-1361600; 6=-1505000; 7=-1808700; 8=-1681700; 9=-1292600; 10=-1139900 }
    94 = @{ 4=-3414800; 5=-1205000; 6=-1343700; 7=-1643400; 8=-1513600; 9=-1187800; 10=-1125800 }
    96 = @{ 4=-554500; 5=-494600; 6=-494500; 7=-500100; 8=-456900; 9=-387600; 10=-358900 }
    100 = @{ 4=2205300; 5=-973800; 6=-1564600; 7=-1301600; 8=-1404400; 9=-649200; 10=169800 }
    101 = @{ 4=89500; 5=-13300; 6=22100; 7=-25700; 8=5600; 9=-2600; 10=-9000 }
    102 = @{ 4=221400; 5=123300; 6=3200; 7=135900; 8=-24800; 9=95300; 10=-10400 }
}

foreach ($row in $changes.Keys) {
    $colmap = $changes[$row]
    foreach ($col in $colmap.Keys) {
        $ws.Cells.Item($row, $col).Value = $colmap[$col]
    }
}